# Insert a new weekly price record as row 37, shifting existing rows
# 37-56 down to 38-57 (matches the "Fruta / hortaliza, semanal" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(37).Insert()

$ws.Range("A37").Value = 1
$ws.Range("B37").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C37").Value = "Arica y Parinacota"
$ws.Range("D37").Value = 44900
$ws.Range("E37").Value = 15
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103004
$ws.Range("J37").Value = "Durazno"
$ws.Range("K37").Value = "Early Majestic"
$ws.Range("L37").Value = "Segunda"
$ws.Range("M37").Value = 200
$ws.Range("N37").Value = 19000
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 19500
$ws.Range("Q37").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R37").Value = "Región de O'Higgins"
$ws.Range("S37").Value = 1083
$ws.Range("T37").Value = 18
